$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.771.50"
$ws.Range("E2").Value = "  +4.74%  "
$ws.Range("D3").Value = "3.101.78"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.097.44"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +9.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "3.614.80"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "3.101.97"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").Value = "62.695.74"
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.727"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.59%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("E33").Value = "  +7.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.27%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0817"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "425.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.15%  "
$ws.Range("D42").Value = "2.911.06"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0368"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.277"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.91%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.56%  "
